$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle")

$ws.Range("E5").Value = "Yes"
$ws.Range("G5").Value = "yes"
$ws.Range("F5").Value = "no"
$ws.Range("H5").Value = "no"

$ws.Range("A6").Select()
